# Added 4wk low sales check: updated Inventory Coverage / Seasonality Index
# figures on the Forecast Comparison sheet (and the dependent 16/8/4-week
# forecast totals on the Summary sheet) to reflect the new low-sales
# detection logic.

$wb = $excel.ActiveWorkbook

$forecast = $wb.Worksheets.Item("Forecast Comparison")
$summary  = $wb.Worksheets.Item("Summary")

# Row -> [Inventory Coverage (H), Seasonality Index (L)]
$rows = @{
    2  = @(18.7,  1.18)
    3  = @(17.7,  1.02)
    4  = @(16.7,  1.09)
    5  = @(15.7,  1)
    6  = @(13.52, 1.19)
    7  = @(12.52, 0.95)
    8  = @(12.52, 1.05)
    9  = @(11.52, 1.2)
    10 = @(9.68,  1.07)
    11 = @(8.68,  1.17)
    12 = @(8.35,  0.89)
    13 = @(8.89,  0.84)
    14 = @(8.82,  0.8)
    15 = @(7.82,  0.99)
    16 = @(6.11,  0.96)
    17 = @(5.71,  1.05)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $forecast.Range("H$r").Value = $vals[0]
    $forecast.Range("L$r").Value = $vals[1]
}

# MyForecast for W22 (row 14) dropped from 2 to 1 units.
$forecast.Range("D14").Value = 1

# Summary sheet: forecast totals shift upward with the revised weekly figures.
# (Leading apostrophe keeps these as text cells, matching the existing
# column formatting where all Value entries are stored as text.)
$summary.Range("B9").Value  = "'35"
$summary.Range("B10").Value = "'19"
$summary.Range("B11").Value = "'9"
